$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date header row (row 1, columns B:L) - shifted forward by one week
$dates = @(45685, 45692, 45699, 45706, 45713, 45720, 45727, 45734, 45741, 45748, 45755)
for ($i = 0; $i -lt $dates.Length; $i++) {
    $ws.Cells.Item(1, 2 + $i).Value = $dates[$i]
}

# Updated forecast values for each row (rows 2-12, columns B:L)
$data = @{
    2  = @(43.4, 43.2, 43.7, 44.8, 46.5, 48.8, 51.6, 54.8, 58.3, 61.9, 65.40000000000001)
    3  = @(42.2, 41.5, 41.4, 41.7, 42.5, 43.8, 45.5, 47.7, 50.2, 53, 56)
    4  = @(39.3, 37.5, 36.4, 36.2, 36.9, 38.6, 41.1, 44.4, 48.2, 52.3, 56.4)
    5  = @(67.09999999999999, 67.5, 67.7, 67.90000000000001, 68.09999999999999, 68.40000000000001, 68.8, 69.5, 70.40000000000001, 71.5, 72.8)
    6  = @(43.1, 43.7, 45.2, 47.3, 50, 53.1, 56.5, 60.1, 63.8, 67.3, 70.59999999999999)
    7  = @(39.1, 37, 35.9, 35.8, 37, 39.3, 42.7, 47.1, 52.1, 57.5, 62.9)
    8  = @(31.4, 29.9, 29.7, 30.9, 33.5, 37.4, 42.5, 48.4, 54.9, 61.5, 67.8)
    9  = @(33.5, 34.5, 36.3, 38.7, 41.8, 45.5, 49.6, 54.3, 59.3, 64.5, 69.7)
    10 = @(48.7, 50, 51.7, 53.8, 56.2, 58.8, 61.5, 64.09999999999999, 66.7, 69.09999999999999, 71.09999999999999)
    11 = @(39.1, 41, 43.5, 46.5, 49.5, 52.6, 55.4, 58, 60.5, 62.7, 64.8)
    12 = @(50, 49.4, 49, 49.1, 49.7, 50.9, 52.5, 54.5, 56.9, 59.4, 61.8)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, 2 + $i).Value = $vals[$i]
    }
}
